{"js": "// Replace the date line and every \"AxB=\" multiplication prompt in the\n// table with the new values from the target revision. Each \"old\" value\n// is unique within the document, so a plain text search+replace per\n// pair is unambiguous and keeps the existing run formatting (font,\n// size, etc.) untouched because insertText(..., replace) only swaps\n// the text inside the matched range.\n\nconst replacements = [\n  [\"2025-05-12 Monday\", \"2025-05-13 Tuesday\"],\n  [\"972\u00d77=\", \"830\u00d76=\"],\n  [\"274\u00d74=\", \"604\u00d76=\"],\n  [\"132\u00d79=\", \"337\u00d73=\"],\n  [\"753\u00d73=\", \"592\u00d72=\"],\n  [\"889\u00d74=\", \"850\u00d74=\"],\n  [\"430\u00d72=\", \"142\u00d79=\"],\n  [\"476\u00d76=\", \"652\u00d75=\"],\n  [\"344\u00d79=\", \"507\u00d78=\"],\n  [\"234\u00d74=\", \"509\u00d76=\"],\n  [\"155\u00d79=\", \"678\u00d74=\"],\n  [\"825\u00d79=\", \"672\u00d77=\"],\n  [\"252\u00d76=\", \"705\u00d79=\"],\n  [\"644\u00d74=\", \"883\u00d74=\"],\n  [\"551\u00d75=\", \"329\u00d72=\"],\n  [\"303\u00d76=\", \"929\u00d73=\"],\n  [\"827\u00d75=\", \"672\u00d73=\"],\n  [\"457\u00d78=\", \"329\u00d72=\"],\n  [\"543\u00d77=\", \"192\u00d79=\"],\n  [\"677\u00d74=\", \"667\u00d78=\"],\n  [\"388\u00d72=\", \"203\u00d78=\"],\n  [\"359\u00d75=\", \"709\u00d74=\"],\n  [\"762\u00d79=\", \"681\u00d79=\"],\n  [\"332\u00d78=\", \"119\u00d73=\"],\n  [\"177\u00d73=\", \"395\u00d78=\"],\n  [\"412\u00d73=\", \"445\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"AxB=\" multiplication prompt in the\n# table with the new values from the target revision. Each \"old\" value\n# is unique within the document, so a plain Find/Replace per pair is\n# unambiguous and the surrounding run formatting (font, size, etc.) is\n# left untouched since Find/Replace only swaps the matched text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-12 Monday\", \"2025-05-13 Tuesday\"),\n    @(\"972\u00d77=\", \"830\u00d76=\"),\n    @(\"274\u00d74=\", \"604\u00d76=\"),\n    @(\"132\u00d79=\", \"337\u00d73=\"),\n    @(\"753\u00d73=\", \"592\u00d72=\"),\n    @(\"889\u00d74=\", \"850\u00d74=\"),\n    @(\"430\u00d72=\", \"142\u00d79=\"),\n    @(\"476\u00d76=\", \"652\u00d75=\"),\n    @(\"344\u00d79=\", \"507\u00d78=\"),\n    @(\"234\u00d74=\", \"509\u00d76=\"),\n    @(\"155\u00d79=\", \"678\u00d74=\"),\n    @(\"825\u00d79=\", \"672\u00d77=\"),\n    @(\"252\u00d76=\", \"705\u00d79=\"),\n    @(\"644\u00d74=\", \"883\u00d74=\"),\n    @(\"551\u00d75=\", \"329\u00d72=\"),\n    @(\"303\u00d76=\", \"929\u00d73=\"),\n    @(\"827\u00d75=\", \"672\u00d73=\"),\n    @(\"457\u00d78=\", \"329\u00d72=\"),\n    @(\"543\u00d77=\", \"192\u00d79=\"),\n    @(\"677\u00d74=\", \"667\u00d78=\"),\n    @(\"388\u00d72=\", \"203\u00d78=\"),\n    @(\"359\u00d75=\", \"709\u00d74=\"),\n    @(\"762\u00d79=\", \"681\u00d79=\"),\n    @(\"332\u00d78=\", \"119\u00d73=\"),\n    @(\"177\u00d73=\", \"395\u00d78=\"),\n    @(\"412\u00d73=\", \"445\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
